$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price cells that receive a new numeric-looking value to Text
# format first, so the new value is stored as a string (matching the
# original t="inlineStr" cell type) instead of being auto-parsed as a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.758.51'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.919.73'
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '242.14'
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4924'
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '0.2996'
$ws.Range("E8").Value = '  +1.82%  '
$ws.Range("D9").Value = '0.06795'
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").Value = '1.915.84'
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").Value = '17.27'
$ws.Range("D12").Value = '0.07351'
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").Value = '5.238'
$ws.Range("E13").Value = '  +3.86%  '
$ws.Range("D14").Value = '89.16'
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("D15").Value = '0.6796'
$ws.Range("E15").Value = '  +0.56%  '
$ws.Range("D16").Value = '30.733.53'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '0.000008017'
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").Value = '13.68'
$ws.Range("E18").Value = '  +4.08%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '2.163.84'
$ws.Range("E20").Value = '  +1.49%  '
$ws.Range("D21").Value = '5.419'
$ws.Range("E21").Value = '  +12.55%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = '202.22'
$ws.Range("E23").Value = '  +10.02%  '
$ws.Range("D24").Value = '6.342'
$ws.Range("D25").Value = '9.714'
$ws.Range("E25").Value = '  +4.15%  '
$ws.Range("D26").Value = '160.97'
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("D27").Value = '18.92'
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("E28").Value = '  +3.87%  '
$ws.Range("D29").Value = '1.472'
$ws.Range("E29").Value = '  +5.24%  '
$ws.Range("D30").Value = '4.377'
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("D31").Value = '0.09188'
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("D32").Value = '4.097'
$ws.Range("E32").Value = '  +2.61%  '
$ws.Range("D33").Value = '0.05309'
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("D34").Value = '0.7495'
$ws.Range("D35").Value = '1.131'
$ws.Range("E35").Value = '  +1.92%  '
$ws.Range("D36").Value = '2.702'
$ws.Range("E36").Value = '  -1.80%  '
$ws.Range("D37").Value = '0.01866'
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("D38").Value = '2.727'
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("D39").Value = '0.9329'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").Value = '2.097'
$ws.Range("E40").Value = '  -1.76%  '
$ws.Range("D41").Value = '0.4517'
$ws.Range("E41").Value = '  +2.44%  '
$ws.Range("D42").Value = '73.19'
$ws.Range("E42").Value = '  +26.88%  '
$ws.Range("D43").Value = '5.986'
$ws.Range("E43").Value = '  +4.49%  '
$ws.Range("D44").Value = '107.67'
$ws.Range("E44").Value = '  +2.37%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '0.1405'
$ws.Range("E45").Value = '  +5.37%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '7.756'
$ws.Range("E47").Value = '  +2.49%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '36.14'
$ws.Range("E48").Value = '  +8.05%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.205'
$ws.Range("E49").Value = '  +6.69%  '
$ws.Range("D50").Value = '0.05942'
$ws.Range("E50").Value = '  +1.62%  '
$ws.Range("D51").Value = '0.4073'
$ws.Range("E51").Value = '  +3.91%  '
